# Auto update Excel log
# Appends newly-logged sensor readings to the "PIR" and "Humidity" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append rows 208-220 (Bathroom, No Motion / Inactive)
# ---------------------------------------------------------------------
$pirWs = $wb.Worksheets.Item("PIR")

$pirTimestamps = @("18:38:29", "18:38:32", "18:38:37", "18:38:42", "18:38:47", "18:38:52", "18:38:57", "18:39:02", "18:39:07", "18:39:12", "18:39:17", "18:39:22", "18:39:28")

$startRow = 208
for ($i = 0; $i -lt $pirTimestamps.Count; $i++) {
    $r = $startRow + $i
    $pirWs.Cells.Item($r, 1).Value = "'2026-01-30"
    $pirWs.Cells.Item($r, 2).Value = $pirTimestamps[$i]
    $pirWs.Cells.Item($r, 3).Value = "18:00"
    $pirWs.Cells.Item($r, 4).Value = "Bathroom"
    $pirWs.Cells.Item($r, 5).Value = "No Motion"
    $pirWs.Cells.Item($r, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 137-144 (Bathroom, % reading / Active)
# ---------------------------------------------------------------------
$humWs = $wb.Worksheets.Item("Humidity")

$humTimestamps = @("18:38:29", "18:38:38", "18:38:48", "18:38:53", "18:39:08", "18:39:13", "18:39:18", "18:39:28")
$humValues = @("85.9%", "85.9%", "85.9%", "85.0%", "86.0%", "86.0%", "86.0%", "86.0%")

$startRow2 = 137
for ($i = 0; $i -lt $humTimestamps.Count; $i++) {
    $r = $startRow2 + $i
    $humWs.Cells.Item($r, 1).Value = "'2026-01-30"
    $humWs.Cells.Item($r, 2).Value = $humTimestamps[$i]
    $humWs.Cells.Item($r, 3).Value = "18:00"
    $humWs.Cells.Item($r, 4).Value = "Bathroom"
    $humWs.Cells.Item($r, 5).Value = "'" + $humValues[$i]
    $humWs.Cells.Item($r, 6).Value = "Active"
}
